# Quotebook App — template update for version 1.1.0
# Adds three new worksheets (QOTD, Reader, AppFlag), a "User Type" column to
# the Users sheet, and hides the now-unused ID column on Users.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Users sheet: add a "User Type" column (D) and hide column C.
# ---------------------------------------------------------------------------
$users = $wb.Worksheets.Item("Users")

$users.Range("D1").Value = "User Type (A = Admin, N = Normal, I = Invisible)"
$users.Range("D1").Font.Size = 9
$users.Range("D1").Font.Color = 2039583

$users.Range("D2").Value = "A"
$users.Range("D3").Value = "N"
$users.Range("D4").Value = "I"
$users.Range("D2:D4").Font.Size = 10
$users.Range("D2:D4").Font.Color = 0

$users.Range("A1:D9").RowHeight = 12.75

$users.Columns.Item(3).EntireColumn.Hidden = $true
$users.Columns.Item(3).ColumnWidth = 0

$users.Range("G9").Select()

# ---------------------------------------------------------------------------
# 2. New sheet: QOTD ("Quote of the day") — reuses the Quotes header row.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$qotd = $wb.Worksheets.Add($null, $lastSheet)
$qotd.Name = "QOTD"

$qotd.Columns.Item(1).ColumnWidth = 10.7109375
$qotd.Columns.Item(2).ColumnWidth = 9.85546875
$qotd.Columns.Item(3).ColumnWidth = 7.85546875
$qotd.Columns.Item(4).ColumnWidth = 61.28515625

$qotd.Range("A1").Value = "Timestamp"
$qotd.Range("B1").Value = "Entered by"
$qotd.Range("C1").Value = "Said by"
$qotd.Range("D1").Value = "Quote"

$qotd.Range("A3").Select()

# ---------------------------------------------------------------------------
# 3. New sheet: Reader — Google-Drive-backed file reference table.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$reader = $wb.Worksheets.Add($null, $lastSheet)
$reader.Name = "Reader"

$reader.Columns.Item(1).ColumnWidth = 25.28515625
$reader.Columns.Item(2).ColumnWidth = 28.140625
$reader.Columns.Item(3).ColumnWidth = 38.140625
$reader.Columns.Item(4).ColumnWidth = 15.7109375

$reader.Range("A1").Value = "File Name"
$reader.Range("B1").Value = "Friendly Name"
$reader.Range("C1").Value = "File ID"
$reader.Range("D1").Value = "File Size (Bytes)"
$reader.Range("A1:D1").Font.Size = 10
$reader.Range("A1:D1").Font.Color = 2039583

$reader.Range("A2").Value = "my_pdf_file"
$reader.Range("B2").Value = "My PDF File"
$reader.Range("C2").Value = "The File ID is the last string of characters in the URL for your Google Drive File. For example, in the URL https://docs.google.com/spreadsheets/d/1qpyC0XzvTcKT6EISywvqESX3A0MwQoFDE8p-Bll4hps/edit#gid=0, the File ID is 1qpyC0XzvTcKT6EISywvqESX3A0MwQoFDE8p-Bll4hps. You can extract the File ID from your Google Drive files in the same way."
$reader.Range("D2").Value = 50000
$reader.Range("A2:B2").Font.Size = 10
$reader.Range("A2:B2").Font.Color = 0

$reader.Range("E5").Select()

# ---------------------------------------------------------------------------
# 4. New sheet: AppFlag — a single cell used as a kill switch for the app.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$appFlag = $wb.Worksheets.Add($null, $lastSheet)
$appFlag.Name = "AppFlag"

$appFlag.Range("B1").Value = "In Cell A1, put a Y to disable access to the app. Whenever there is a Y there, users will not be able to log in. "
$appFlag.Range("B1").Font.Size = 10
$appFlag.Range("B1").Font.Color = 0

$appFlag.Range("F21").Select()
